$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at row 167; everything currently at/after row 167
# (old rows 167-179) shifts down to rows 168-180.
$ws.Rows.Item(167).Insert()

# Populate the newly inserted row 167 with the new record.
$ws.Range("A167").Value = 5
$ws.Range("B167").Value = 'Macroferia Regional de Talca'
$ws.Range("C167").Value = 'Maule'
$ws.Range("D167").Value = 44714
$ws.Range("E167").Value = 7
$ws.Range("F167").Value = 100112017
$ws.Range("G167").Value = 'Apio'
$ws.Range("H167").Value = 'Americana (o)'
$ws.Range("I167").Value = 'Primera'
$ws.Range("J167").Value = 600
$ws.Range("K167").Value = 6000
$ws.Range("L167").Value = 6000
$ws.Range("M167").Value = 6000
$ws.Range("N167").Value = '$/docena de matas'
$ws.Range("O167").Value = 'Provincia del Elquí'
$ws.Range("P167").Value = 1000
$ws.Range("Q167").Value = 6
$ws.Range("R167").Value = 'Hortaliza'

# Make sure the date cell keeps the same date number format as the rest of
# column D (style carried over by Insert, but set explicitly to be safe).
$ws.Range("D167").NumberFormat = $ws.Range("D168").NumberFormat
